$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (report number + date range)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value  = "Volume 30   Number  14"
$ws.Range("C9").Value  = "Report Covering the Week  4/3/2023  Through  4/9/2023"

# ---------------------------------------------------------------------------
# Helper behaviour notes:
#  - Reading a COM property in this host requires explicit call syntax, e.g.
#    $range.Value() rather than bare $range.Value.
#  - To flip a cell from the "text placeholder" convention (shared string
#    "0", General-format style) to a real number, or back, we first copy the
#    destination (number) formatting from a sibling cell that already has
#    the wanted style via PasteSpecial(xlPasteFormats = -4122); this reuses
#    the existing style index instead of minting a brand new one.
# ---------------------------------------------------------------------------

# Row 16 (Robbery) - pct-change figures refreshed
$ws.Range("M16").Value = -66.666666666666
$ws.Range("N16").Value = -92.857142857142

# Row 17 (Fel. Assault) - C17 goes from text "0" placeholder to a real count
$c17 = $ws.Range("C17")
$c17.NumberFormat = "#,##0"
$c17.Value = 1
$ws.Range("I17").Value = 3
$ws.Range("K17").Value = 50
$ws.Range("L17").Value = 200
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = -75

# Row 19 (Gr. Larceny) - C19 goes from a real count back to the text "0" placeholder
$dst19 = $ws.Range("C19")
$fmtSrc19 = $ws.Range("D19")
$fmtSrc19.Copy()
$dst19.PasteSpecial(-4122)
$dst19.Value = "'0"
$fmtSrc19.Copy()
$dst19.PasteSpecial(-4122)
$ws.Range("M19").Value = -44.444444444444
$ws.Range("N19").Value = -72.222222222222

# Row 21 (TOTAL)
$ws.Range("I21").Value = 11
$ws.Range("K21").Value = 10
$ws.Range("L21").Value = 57.142857142857
$ws.Range("M21").Value = -31.25
$ws.Range("N21").Value = -84.057971014492

# Row 24 (Petit Larceny)
$ws.Range("F24").Value = 2
$ws.Range("H24").Value = 100
$ws.Range("L24").Value = 80
$ws.Range("M24").Value = 125

# Row 25 (Misd. Assault) - C25 goes from a real count back to the text "0" placeholder
$dst25 = $ws.Range("C25")
$fmtSrc25 = $ws.Range("D25")
$fmtSrc25.Copy()
$dst25.PasteSpecial(-4122)
$dst25.Value = "'0"
$fmtSrc25.Copy()
$dst25.PasteSpecial(-4122)
$ws.Range("L25").Value = 20
$ws.Range("M25").Value = -40

# Row 27 (Other Sex Crimes)
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 4
$ws.Range("I27").Value = 5
$ws.Range("K27").Value = 400
